# update RMP data source
# Adds a new row (23) to the indicator search table describing the
# EPA Risk Management Plan (RMP) facilities dataset, and nudges a few
# sheet-level display settings to match the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row -----------------------------------------------------
$ws.Range("A23").Value = "Risk Management Plan (RMP)"
$ws.Range("B23").Value = "CONUS, AK, HI"
$ws.Range("C23").Value = "Y"
$ws.Range("D23").Value = ".csv"
$ws.Range("E23").Value = "points"
$ws.Range("F23").Value = "1.8MB"
$ws.Range("G23").Value = "EPA"
$ws.Range("H23").Value = "https://hifld-geoplatform.opendata.arcgis.com/datasets/geoplatform::epa-emergency-response-er-risk-management-plan-rmp-facilities/explore?location=35.878600%2C-113.806709%2C3.92"
$ws.Range("I23").Value = "Last Updated April 2022"

# The link text in H23 was pasted in as plain (non-hyperlinked) text
# carrying a web-page style font: 13pt Arial, dark gray.
$h23 = $ws.Range("H23")
$h23.Font.Name = "Arial"
$h23.Font.Size = 13
$h23.Font.Color = 1776411

# Row 23 renders a bit taller to fit that larger font.
$ws.Rows.Item(23).RowHeight = 16.5

# --- Selection / view --------------------------------------------------
$ws.Range("I24").Select() | Out-Null

# --- Page setup ---------------------------------------------------------
$ws.PageSetup.Orientation = 1
